$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "year" column (A) currently stores a date serial (e.g. 43221, 43374).
# It should instead hold the plain text value "2018" for both data rows.
# Using a leading apostrophe forces Excel to store it as text (so it lands
# in the shared-string table) rather than re-interpreting it as a number.
$ws.Range("A2").Value = "'2018"
$ws.Range("A3").Value = "'2018"

# Drop the inherited date number-format/style from those two cells so they
# fall back to the sheet's default (General) style, matching a plain text
# cell with no special formatting.
$ws.Range("A2").ClearFormats()
$ws.Range("A3").ClearFormats()
